# Updated symbol list on Wed Feb 15 18:53:02 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) crypto
# quotes on the active sheet. Source data is written as literal text
# (the workbook stores these columns as plain strings, e.g. "302.98" /
# "1.98%", not numbers/percentages), so each value is entered with a
# leading apostrophe to force Excel's text interpretation instead of
# letting it auto-convert to a number or a percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.98"
$ws.Range("E2").Value = "'1.98%"
$ws.Range("D3").Value = "'44.13"
$ws.Range("E3").Value = "'6.75%"
$ws.Range("D4").Value = "'5.107"
$ws.Range("E4").Value = "'1.94%"
$ws.Range("D5").Value = "'0.07736"
$ws.Range("E5").Value = "'3.10%"
$ws.Range("E6").Value = "'1.16%"
$ws.Range("D7").Value = "'1.617"
$ws.Range("E7").Value = "'2.65%"
$ws.Range("E8").Value = "'13.09%"
$ws.Range("E9").Value = "'5.57%"
$ws.Range("D10").Value = "'0.1865"
$ws.Range("E10").Value = "'1.48%"
$ws.Range("D11").Value = "'0.09271"
$ws.Range("E11").Value = "'4.68%"
$ws.Range("D12").Value = "'0.04149"
$ws.Range("E12").Value = "'0.89%"
$ws.Range("D13").Value = "'0.1049"
$ws.Range("E13").Value = "'-0.45%"
$ws.Range("D14").Value = "'0.001281"
$ws.Range("E14").Value = "'-0.12%"
$ws.Range("D15").Value = "'0.005762"
$ws.Range("E15").Value = "'-0.27%"
$ws.Range("D17").Value = "'3.351"
$ws.Range("E17").Value = "'0.17%"
$ws.Range("E18").Value = "'-2.94%"
$ws.Range("D19").Value = "'0.3352"
$ws.Range("E19").Value = "'1.04%"
$ws.Range("D20").Value = "'8.038"
$ws.Range("E20").Value = "'0.70%"
$ws.Range("D21").Value = "'0.1376"
$ws.Range("E21").Value = "'-3.06%"
$ws.Range("E22").Value = "'7.10%"
$ws.Range("D23").Value = "'0.04189"
$ws.Range("E23").Value = "'3.38%"
$ws.Range("D24").Value = "'0.001279"
$ws.Range("E24").Value = "'1.14%"
$ws.Range("D25").Value = "'0.004408"
$ws.Range("E25").Value = "'13.64%"
$ws.Range("E26").Value = "'9.59%"
$ws.Range("D38").Value = "'0.02512"
$ws.Range("E38").Value = "'4.04%"
$ws.Range("D39").Value = "'0.05311"
$ws.Range("E39").Value = "'1.79%"
$ws.Range("D40").Value = "'0.005826"
$ws.Range("E40").Value = "'-2.65%"
$ws.Range("D41").Value = "'0.007732"
$ws.Range("E41").Value = "'-0.86%"
$ws.Range("D42").Value = "'0.1359"
$ws.Range("E42").Value = "'2.55%"
$ws.Range("D43").Value = "'0.007356"
$ws.Range("E43").Value = "'-0.16%"
$ws.Range("D44").Value = "'0.007507"
$ws.Range("E44").Value = "'-7.48%"
$ws.Range("D45").Value = "'0.3018"
$ws.Range("E45").Value = "'1.64%"
$ws.Range("D46").Value = "'0.00006681"
$ws.Range("E46").Value = "'6.45%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("D48").Value = "'0.06430"
$ws.Range("E48").Value = "'39.14%"
$ws.Range("D49").Value = "'0.00002095"
$ws.Range("E49").Value = "'-0.15%"
$ws.Range("E50").Value = "'-0.15%"
